$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("I3").Comment
$c.Delete()
$ws.Range("H2").AddCommentThreaded("BUSCAR código de país ISO 3166-1 ")
